$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "QNCQC5"
$ws.Range("B13").Value = "Calibrador Vernier"
$ws.Range("C13").Value = "150mm."
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 130000
$ws.Range("F13").Value = 4
$ws.Range("G13").Value = 4
$ws.Range("H13").Formula = "=(E13-D13)*G13"
$ws.Range("I13").Formula = "=D13*F13"
$ws.Range("J13").Value = 0
